# Generate Report for Handoff
# Marks the 9830e1ca-... and e5952e6a-... files as "Ready for handoff"
# across the Overview / zh-cn / de-de sheets, bumps their handoff
# timestamps, widens the Error Detail column, and records a
# "handback file is not the latest" warning for each of them.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("E4").Value = "Ready for handoff"
$ws.Range("F4").Value = "Ready for handoff"
$ws.Range("G4").Value = "2016-08-18 08:27:21"

$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"
$ws.Range("G5").Value = "2016-08-18 08:27:21"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

# widen the "Error Detail" column (P) to fit the new warning text
$ws.Columns.Item(16).ColumnWidth = 39.17

$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("H4").Value = "2016-08-18 08:27:15"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fac60ff742a93a773685578f019c71c237b5acd6/e2e/9830e1ca-4d6e-413f-a8e7-043487110af3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d6f159da09b1f5c221736d6d067b641d06723c1/e2e/9830e1ca-4d6e-413f-a8e7-043487110af3.md."

$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("H5").Value = "2016-08-18 08:27:15"
$ws.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fac60ff742a93a773685578f019c71c237b5acd6/e2e/e5952e6a-8cd2-443b-80e4-a21aff8e3ef9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d6f159da09b1f5c221736d6d067b641d06723c1/e2e/e5952e6a-8cd2-443b-80e4-a21aff8e3ef9.md."

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

# widen the "Error Detail" column (P) to fit the new warning text
$ws.Columns.Item(16).ColumnWidth = 39.17

$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("H4").Value = "2016-08-18 08:27:21"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fac60ff742a93a773685578f019c71c237b5acd6/e2e/9830e1ca-4d6e-413f-a8e7-043487110af3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d6f159da09b1f5c221736d6d067b641d06723c1/e2e/9830e1ca-4d6e-413f-a8e7-043487110af3.md."

$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("H5").Value = "2016-08-18 08:27:21"
$ws.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fac60ff742a93a773685578f019c71c237b5acd6/e2e/e5952e6a-8cd2-443b-80e4-a21aff8e3ef9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d6f159da09b1f5c221736d6d067b641d06723c1/e2e/e5952e6a-8cd2-443b-80e4-a21aff8e3ef9.md."
